$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Handled spinner" -> flip RunToTest flag between the existing rows 21 and 23
$ws.Range("A21").Value = "N"
$ws.Range("A23").Value = "Y"

# "and Billet To updated ID" -> append six more consignment rows (24-29) with
# new/updated ConsignmentID values, repeating the same alternating pattern
# used by rows 2-23.
$newRows = @(
    @("N", "UAT42092143", "CELW01", "Cell", "S", "Surekha", 123, "aaaaaa", "qqqqqqq", 1, "SP12345678"),
    @("N", "UAT42092145", "CELW01", "Cell", "A", "Anand",   456, "wwww",   "eeeeeee", 1, "SP12345678"),
    @("N", "UAT42092146", "CELW01", "Cell", "S", "Surekha", 123, "aaaaaa", "qqqqqqq", 1, "SP12345678"),
    @("N", "UAT42092147", "CELW01", "Cell", "A", "Anand",   456, "wwww",   "eeeeeee", 1, "SP12345678"),
    @("N", "UAT42092148", "CELW01", "Cell", "S", "Surekha", 123, "aaaaaa", "qqqqqqq", 1, "SP12345678"),
    @("N", "UAT42092149", "CELW01", "Cell", "A", "Anand",   456, "wwww",   "eeeeeee", 1, "SP12345678")
)

$rowIndex = 24
foreach ($rowData in $newRows) {
    $ws.Cells.Item($rowIndex, 1).Value  = $rowData[0]
    $ws.Cells.Item($rowIndex, 2).Value  = $rowData[1]
    $ws.Cells.Item($rowIndex, 3).Value  = $rowData[2]
    $ws.Cells.Item($rowIndex, 4).Value  = $rowData[3]
    $ws.Cells.Item($rowIndex, 5).Value  = $rowData[4]
    $ws.Cells.Item($rowIndex, 6).Value  = $rowData[5]
    $ws.Cells.Item($rowIndex, 7).Value  = $rowData[6]
    $ws.Cells.Item($rowIndex, 8).Value  = $rowData[7]
    $ws.Cells.Item($rowIndex, 9).Value  = $rowData[8]
    $ws.Cells.Item($rowIndex, 10).Value = $rowData[9]
    $ws.Cells.Item($rowIndex, 11).Value = $rowData[10]
    $rowIndex++
}

# Restore the scrolled viewport / active selection left by the user after
# keying in the extra rows.
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("B35").Select() | Out-Null
